$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Trends Status sheet - update values
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("B2").Value = 5
$wsTrends.Range("C2").Value = 3
$wsTrends.Range("D2").Value = 50
$wsTrends.Range("E2").Value = 15.8

$wsTrends.Range("B3").Value = 1
$wsTrends.Range("C3").Value = 7
$wsTrends.Range("D3").Value = 10
$wsTrends.Range("E3").Value = 36.8

$wsTrends.Range("B4").Value = 3
$wsTrends.Range("C4").Value = 8
$wsTrends.Range("D4").Value = 30
$wsTrends.Range("E4").Value = 42.1

$wsTrends.Range("B5").Value = 1
$wsTrends.Range("C5").Value = 1
$wsTrends.Range("D5").Value = 10
$wsTrends.Range("E5").Value = 5.3

$wsTrends.Range("B6").Value = 0
$wsTrends.Range("C6").Value = 0
$wsTrends.Range("D6").Value = 0
$wsTrends.Range("E6").Value = 0

$wsTrends.Range("B7").Value = 23
$wsTrends.Range("C7").Value = 30

$wsTrends.Range("B8").Value = 85
$wsTrends.Range("C8").Value = 69

# ---------------------------------------------------------------------------
# 2. Priority Status sheet - update values
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3. Species qualification sheet - update text & values
# ---------------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("A2").Value = "SoIB Assessment"
$wsSpecies.Range("B3").Value = 33
$wsSpecies.Range("C3").Value = 10
$wsSpecies.Range("C4").Value = 19

# ---------------------------------------------------------------------------
# 4. Rename "High Priority break-up" to "Interannual update - High Pri".
#    Its header row (A1:E1) is kept as-is.
# ---------------------------------------------------------------------------
$wsInterannual = $wb.Worksheets.Item("High Priority break-up")
$wsInterannual.Name = "Interannual update - High Pri"

# ---------------------------------------------------------------------------
# 5. Create the new "Major update - High Priority " sheet right after
#    "Interannual update - High Pri", and populate it with the content that
#    used to live in the "High Priority break-up" sheet.
# ---------------------------------------------------------------------------
$wsMajor = $wb.Worksheets.Add([System.Type]::Missing, $wsInterannual)
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108

$wsMajor.Range("A2").Value = "Trend New"
$wsMajor.Range("B2").Value = 12
$wsMajor.Range("C2").Value = 48
$wsMajor.Range("D2").Value = 12
$wsMajor.Range("E2").Value = 48

$wsMajor.Range("A3").Value = "IUCN"
$wsMajor.Range("B3").Value = 13
$wsMajor.Range("C3").Value = 52
$wsMajor.Range("D3").Value = 13
$wsMajor.Range("E3").Value = 52

# ---------------------------------------------------------------------------
# 6. Replace the contents of "Interannual update - High Pri" with the
#    updated data (an extra "Trend Different" row is introduced).
# ---------------------------------------------------------------------------
$wsInterannual.Range("A2").Value = "Trend New"
$wsInterannual.Range("B2").Value = 86
$wsInterannual.Range("C2").Value = 83.5
$wsInterannual.Range("D2").Value = 86
$wsInterannual.Range("E2").Value = 97.7

$wsInterannual.Range("A3").Value = "Trend Different"
$wsInterannual.Range("B3").Value = 3
$wsInterannual.Range("C3").Value = 2.9
$wsInterannual.Range("D3").ClearContents()
$wsInterannual.Range("E3").ClearContents()

$wsInterannual.Range("A4").Value = "IUCN"
$wsInterannual.Range("B4").Value = 14
$wsInterannual.Range("C4").Value = 13.6
$wsInterannual.Range("D4").Value = 2
$wsInterannual.Range("E4").Value = 2.3
